# dias-salvos-2025-04.xlsx — add three more daily rows (2025-04-12..14)
# plus a refreshed "Média" (average) summary row, leaving a gap of blank
# rows between the daily log and the summary row (mirrors the author's
# "juntei todos os arquivos" consolidation edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Push the old "Média" row (row 3) down to make room -------------------
# Inserting 6 rows above it moves it from row 3 to row 9, matching the
# final layout (daily rows end at 11, blank rows 3-8, summary at row 12).
$ws.Range("A3:A8").EntireRow.Insert()

# The insert carried the old "Média"/29 row along with it into row 9 — wipe
# it so we can rebuild rows 9-12 with the new data from scratch.
$ws.Range("A9:C9").ClearContents()

# --- Row 9: 2025-04-12 | 24 | - --------------------------------------------
$ws.Cells.Item(9, 1).NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = "2025-04-12"
$ws.Cells.Item(9, 1).ClearFormats()
$ws.Cells.Item(9, 2).Value = 24
$ws.Cells.Item(9, 3).Value = "-"

# --- Row 10: 2025-04-13 | 25 | - -------------------------------------------
$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = "2025-04-13"
$ws.Cells.Item(10, 1).ClearFormats()
$ws.Cells.Item(10, 2).Value = 25
$ws.Cells.Item(10, 3).Value = "-"

# --- Row 11: 2025-04-14 | 27 | - -------------------------------------------
$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = "2025-04-14"
$ws.Cells.Item(11, 1).ClearFormats()
$ws.Cells.Item(11, 2).Value = 27
$ws.Cells.Item(11, 3).Value = "-"

# --- Row 12: Média | 26 (no "Acima da Meta" column here, same as before) --
$ws.Cells.Item(12, 1).Value = "Média"
$ws.Cells.Item(12, 2).Value = 26
